$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.972.16"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").Value = "2.302.57"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.85"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.92"
$ws.Range("E6").Value = "  -1.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  -1.73%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.76"
$ws.Range("E10").Value = "  -3.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0907"
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.47"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.976"
$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.36"
$ws.Range("E15").Value = "  -4.13%  "
$ws.Range("D16").Value = "2.650.63"
$ws.Range("E16").Value = "  -2.16%  "
$ws.Range("D17").Value = "2.289.34"
$ws.Range("E17").Value = "  -2.42%  "
$ws.Range("D18").Value = "42.063.48"
$ws.Range("E18").Value = "  -1.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.72"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.16"
$ws.Range("E21").Value = "  -5.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.58"
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "259.89"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.96"
$ws.Range("E25").Value = "  +4.27%  "
$ws.Range("E26").Value = "  +0.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.92"
$ws.Range("E27").Value = "  -4.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.84"
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.27"
$ws.Range("E29").Value = "  +2.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.79"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.02"
$ws.Range("E31").Value = "  -6.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0884"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("E33").Value = "  -2.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.88"
$ws.Range("E34").Value = "  -4.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.121"
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("E38").Value = "  +8.89%  "
$ws.Range("E39").Value = "  -2.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.64"
$ws.Range("E40").Value = "  -4.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "103.33"
$ws.Range("E41").Value = "  +16.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.48"
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.88"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.226"
$ws.Range("E44").Value = "  -2.69%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "115.23"
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.14"
$ws.Range("E47").Value = "  +1.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "78.07"
$ws.Range("E48").Value = "  +5.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.05"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("E50").Value = "  -3.32%  "
$ws.Range("E51").Value = "  +2.03%  "
